$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Append a new row to the eCoaching Log Runbook revision-history table
#    (the 3-column table whose last existing row is the 08/15/2022 / TFS 25205
#    "New Submission" entry authored by Lili Huang).
# ---------------------------------------------------------------------------
$logTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $lastRow = $t.Rows.Item($t.Rows.Count)
    if ($lastRow.Cells.Count -eq 3 -and `
        $lastRow.Cells.Item(1).Range.Text -like "08/15/2022*" -and `
        $lastRow.Cells.Item(2).Range.Text -like "*TFS 25205*" -and `
        $lastRow.Cells.Item(3).Range.Text -like "*Lili Huang*") {
        $logTable = $t
        break
    }
}

$newRow = $logTable.Rows.Add()

# Column 1: date
$dateCell = $newRow.Cells.Item(1)
$dateCell.Range.Text = "08/24/2022"

# Column 2: description - two paragraphs
$descCell = $newRow.Cells.Item(2)
$descCell.Range.Text = "TFS 25205;" + [char]13 + "Included web scan fix " + [char]0x2013 + " moment.js version update"

# Column 3: author
$authorCell = $newRow.Cells.Item(3)
$authorCell.Range.Text = "Lili Huang"

# ---------------------------------------------------------------------------
# 2) Bump the Changeset number referenced in the publish steps from 51847 to
#    51867 (TFS 25205 moment.js web-scan fix shipped in a later changeset).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("847", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "867", 2) | Out-Null
